# Hindalco price sheet update (2025-12-20): a new top row is published each
# day, so this commit inserts one new row above the current row 2 (pushing
# every existing row down by one) and fills it with the day's figures.
# Because the sheet keeps a hyperlink object per "Circular Link" cell
# (F2:F131 before the edit), and this engine's row-insert does not carry
# hyperlink anchors along with the shifted cells, the hyperlinks are
# rebuilt from scratch afterward so every F cell's link matches its text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at position 2; this shifts rows 2:192 down to
#    3:193 (values/text move correctly; dimension grows to F193).
$ws.Rows.Item(2).Insert()

# 2) The fresh row has no explicit style yet - copy the (correct) format
#    from the row right below it (old row 2, now row 3) before writing values.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# 3) Populate the new row with today's (2025-12-20) figures - same
#    description/grade/price/circular as the prior top row, new date.
$ws.Range("A2").Value = "20-12-2025"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 304.5
$ws.Range("E2").Value = "17.12.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-17-december-2025.pdf"

# 4) Existing hyperlink objects stayed anchored to their old row numbers
#    instead of following the shifted cells, so every one is now
#    misaligned with the text that moved under it. Wipe them all ...
$ws.Range("F2").Hyperlinks.Delete()

# 5) ... and rebuild one per row for F2:F132 (the new hyperlink span),
#    pointing each cell at the URL that is already sitting in its own text.
for ($r = 2; $r -le 132; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Text
    $ws.Hyperlinks.Add($cell, $url, "", "", $url)
}

# 6) Hyperlinks.Add stamps cells with the built-in blue/underlined
#    "Hyperlink" style; restore the plain data-row style (same as the
#    neighbouring Circular-Date column) so formatting matches the source.
$ws.Range("E2:E132").Copy()
$ws.Range("F2:F132").PasteSpecial(-4122)
